$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 - "text4" label plus the same formula pattern as the other rows.
$ws.Range("A12").Value = "text4"

$ws.Range("B12").Value = 0.55
$ws.Range("C12").Value = 0.6
$ws.Range("D12").Value = 0.9

$ws.Range("E12").Formula = '=IF(($E$1/100*B12/16)+C12<D12,D12,IF(($E$1/100*B12/16)+C12>J12,J12,($E$1/100*B12/16)+C12))'
$ws.Range("F12").Formula = '=IF(($F$1/100*B12/16)+C12<D12,D12,IF(($F$1/100*B12/16)+C12>J12,J12,($F$1/100*B12/16)+C12))'
$ws.Range("G12").Formula = '=IF(($G$1/100*B12/16)+C12<D12,D12,IF(($G$1/100*B12/16)+C12>J12,J12,($G$1/100*B12/16)+C12))'
$ws.Range("H12").Formula = '=IF(($H$1/100*B12/16)+C12<D12,D12,IF(($H$1/100*B12/16)+C12>J12,J12,($H$1/100*B12/16)+C12))'
$ws.Range("I12").Formula = '=IF(($I$1/100*B12/16)+C12<D12,D12,IF(($I$1/100*B12/16)+C12>J12,J12,($I$1/100*B12/16)+C12))'

$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 1

$ws.Range("L12").Formula = '=E12*16'
$ws.Range("M12").Formula = '=F12*16'
$ws.Range("N12").Formula = '=G12*16'
$ws.Range("O12").Formula = '=H12*16'
$ws.Range("P12").Formula = '=I12*16'

$ws.Range("S12").Formula = '="font-size: clamp("&D12&"rem, "&B12&"vw + "&C12&"rem, "&J12&"rem);"'

# Update the active selection to match the recorded cursor position.
$ws.Range("S12").Select()

$wb.Save()
